# Regenerate orders with updated distance/size codes.
# Mapping:
#   Distance: D51 -> D55, D64 -> D69, D80 -> D86
#   Size:     S30 -> S31
# These substrings appear inside Condition, Filename_Left, Filename_Right,
# Distance and Size columns (B, D, E, H, J) for every trial row. The header
# row (row 1) and other columns (Trial, Duration_Seconds, Is_Repeat, Block,
# Face, ConditionID) never contain these tokens, so a blanket text
# substitution over every used cell is safe and idempotent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$nRows = $used.Rows.Count
$nCols = $used.Columns.Count

for ($r = 2; $r -le $nRows; $r++) {
    for ($c = 1; $c -le $nCols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $txt = $cell.Text
        if ($txt -and ($txt -like "*D51*" -or $txt -like "*D64*" -or $txt -like "*D80*" -or $txt -like "*S30*")) {
            $new = $txt -replace "D51", "D55" -replace "D64", "D69" -replace "D80", "D86" -replace "S30", "S31"
            $cell.Value = $new
        }
    }
}
